# GSoC -> Google Summer of Code
# The abbreviation "GSoC" (a single bold run) is expanded in place to the
# full phrase "Google Summer of Code", keeping the existing (bold,
# Garamond, 22pt) character formatting intact.

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "GSoC",                     # FindText
    $true,                      # MatchCase
    $true,                      # MatchWholeWord
    $false,                     # MatchWildcards
    $false,                     # MatchSoundsLike
    $false,                     # MatchAllWordForms
    $true,                      # Forward
    1,                          # Wrap (wdFindContinue)
    $false,                     # Format
    "Google Summer of Code",    # ReplaceWith
    2                           # Replace (wdReplaceAll)
)
